$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I0 / IF headers in I1:J1, matching the style of existing headers (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J33
$values = @(
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(12, 12),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(10, 10),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(4, 5),
    @(8, 8),
    @(6, 6),
    @(4, 4),
    @(7, 7),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
